{"js": "// Replace the worksheet date and each of the 25 division-problem cells\n// with their new values. Every \"from\" value below is unique within the\n// document, so a simple exact-text search + replace for each pair is\n// enough to retarget every cell without disturbing run formatting.\nconst replacements = [\n  { from: \"2025-11-30 Sunday\", to: \"2025-12-01 Monday\" },\n  { from: \"814\u00f76=135, 4\", to: \"155\u00f73=51, 2\" },\n  { from: \"411\u00f79=45, 6\", to: \"345\u00f75=69, 0\" },\n  { from: \"470\u00f76=78, 2\", to: \"473\u00f73=157, 2\" },\n  { from: \"542\u00f76=90, 2\", to: \"365\u00f76=60, 5\" },\n  { from: \"302\u00f75=60, 2\", to: \"972\u00f78=121, 4\" },\n  { from: \"576\u00f76=96, 0\", to: \"253\u00f73=84, 1\" },\n  { from: \"130\u00f73=43, 1\", to: \"489\u00f74=122, 1\" },\n  { from: \"459\u00f77=65, 4\", to: \"653\u00f78=81, 5\" },\n  { from: \"125\u00f74=31, 1\", to: \"896\u00f77=128, 0\" },\n  { from: \"782\u00f78=97, 6\", to: \"186\u00f75=37, 1\" },\n  { from: \"493\u00f79=54, 7\", to: \"692\u00f74=173, 0\" },\n  { from: \"895\u00f75=179, 0\", to: \"416\u00f75=83, 1\" },\n  { from: \"162\u00f72=81, 0\", to: \"420\u00f78=52, 4\" },\n  { from: \"703\u00f76=117, 1\", to: \"278\u00f79=30, 8\" },\n  { from: \"833\u00f77=119, 0\", to: \"884\u00f75=176, 4\" },\n  { from: \"541\u00f78=67, 5\", to: \"978\u00f76=163, 0\" },\n  { from: \"627\u00f73=209, 0\", to: \"178\u00f79=19, 7\" },\n  { from: \"992\u00f72=496, 0\", to: \"750\u00f77=107, 1\" },\n  { from: \"845\u00f73=281, 2\", to: \"351\u00f72=175, 1\" },\n  { from: \"562\u00f76=93, 4\", to: \"868\u00f76=144, 4\" },\n  { from: \"859\u00f75=171, 4\", to: \"800\u00f76=133, 2\" },\n  { from: \"710\u00f74=177, 2\", to: \"354\u00f74=88, 2\" },\n  { from: \"320\u00f76=53, 2\", to: \"581\u00f79=64, 5\" },\n  { from: \"644\u00f72=322, 0\", to: \"269\u00f79=29, 8\" },\n  { from: \"940\u00f72=470, 0\", to: \"462\u00f75=92, 2\" },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${from}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and all 25 division-problem answers with\n# their new values. Every \"find\" value below is unique in the document,\n# so a simple Find/Replace (wdReplaceAll) pass per pair is sufficient\n# and keeps each run's original formatting untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-30 Sunday\", \"2025-12-01 Monday\"),\n    @(\"814\u00f76=135, 4\", \"155\u00f73=51, 2\"),\n    @(\"411\u00f79=45, 6\", \"345\u00f75=69, 0\"),\n    @(\"470\u00f76=78, 2\", \"473\u00f73=157, 2\"),\n    @(\"542\u00f76=90, 2\", \"365\u00f76=60, 5\"),\n    @(\"302\u00f75=60, 2\", \"972\u00f78=121, 4\"),\n    @(\"576\u00f76=96, 0\", \"253\u00f73=84, 1\"),\n    @(\"130\u00f73=43, 1\", \"489\u00f74=122, 1\"),\n    @(\"459\u00f77=65, 4\", \"653\u00f78=81, 5\"),\n    @(\"125\u00f74=31, 1\", \"896\u00f77=128, 0\"),\n    @(\"782\u00f78=97, 6\", \"186\u00f75=37, 1\"),\n    @(\"493\u00f79=54, 7\", \"692\u00f74=173, 0\"),\n    @(\"895\u00f75=179, 0\", \"416\u00f75=83, 1\"),\n    @(\"162\u00f72=81, 0\", \"420\u00f78=52, 4\"),\n    @(\"703\u00f76=117, 1\", \"278\u00f79=30, 8\"),\n    @(\"833\u00f77=119, 0\", \"884\u00f75=176, 4\"),\n    @(\"541\u00f78=67, 5\", \"978\u00f76=163, 0\"),\n    @(\"627\u00f73=209, 0\", \"178\u00f79=19, 7\"),\n    @(\"992\u00f72=496, 0\", \"750\u00f77=107, 1\"),\n    @(\"845\u00f73=281, 2\", \"351\u00f72=175, 1\"),\n    @(\"562\u00f76=93, 4\", \"868\u00f76=144, 4\"),\n    @(\"859\u00f75=171, 4\", \"800\u00f76=133, 2\"),\n    @(\"710\u00f74=177, 2\", \"354\u00f74=88, 2\"),\n    @(\"320\u00f76=53, 2\", \"581\u00f79=64, 5\"),\n    @(\"644\u00f72=322, 0\", \"269\u00f79=29, 8\"),\n    @(\"940\u00f72=470, 0\", \"462\u00f75=92, 2\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $findText\"\n    }\n}\n"}
